$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.998.53'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.450.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0976'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("E10").Value = '  -1.98%  '
$ws.Range("E11").Value = '  -4.21%  '
$ws.Range("E12").Value = '  -4.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.885.06'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '57.917.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.77%  '
$ws.Range("E16").Value = '  -2.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.451.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.44%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.403'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.564.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.51%  '
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '174.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.66%  '
$ws.Range("E30").Value = '  -2.95%  '
$ws.Range("E31").Value = '  -2.44%  '
$ws.Range("E32").Value = '  -3.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.28%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  -2.81%  '
$ws.Range("E37").Value = '  -7.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.30'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.809'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("E41").Value = '  -3.91%  '
$ws.Range("E42").Value = '  -3.13%  '
$ws.Range("E43").Value = '  -3.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '125.39'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '258.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.04%  '
$ws.Range("E46").Value = '  -5.42%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E49").Value = '  -2.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.02'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.93%  '
$ws.Range("E51").Value = '  -5.04%  '
